# Expense Sheet -- in Progress
# Insert a new blank separator row above the "User Administration" block
# (shifts everything from the old row 13 down by one row) and append two
# new feature blocks ("Loggin & Instrumentation" and "File Manager") at
# the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the "User Administration" block (old rows 13-24) down by one row.
$ws.Rows("13:13").Insert()

# New block: Loggin & Instrumentation (rows 28-30)
$ws.Range("B28").Value = "Loggin & Instrumentation"
$ws.Range("C28").Value = "Using Serilog with .net Core projects for DI"
$ws.Range("C29").Value = "Add Serilog Logger to Aspnet Core Web"
$ws.Range("C30").Value = "Add Serilog Logger to Web Apis'"

# New block: File Manager (row 33)
$ws.Range("B33").Value = "File Manager"
$ws.Range("C33").Value = "Add files grouped by folders… for salary slips…. Look at UI in theme `"File Manager`""

# Column B needs to widen to fit the new "Loggin & Instrumentation" text.
$ws.Columns("B:B").ColumnWidth = 20.17

# Move the selection/view down to where the new rows were added.
$ws.Range("C34").Select() | Out-Null
